$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$defaultStyleRef = $ws.Range("B2")

$ws.Range('D2').Value = '44.039.15'
$ws.Range('E2').Value = '  -1.15%  '
$ws.Range('D3').Value = '2.359.99'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '239.82'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = $defaultStyleRef.Style
$ws.Range('E6').Value = '  +0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.07'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = $defaultStyleRef.Style
$ws.Range('E7').Value = '  +1.00%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.596'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = $defaultStyleRef.Style
$ws.Range('E9').Value = '  +7.75%  '
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.28'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = $defaultStyleRef.Style
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '32.27'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = $defaultStyleRef.Style
$ws.Range('E12').Value = '  +7.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.28'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = $defaultStyleRef.Style
$ws.Range('E13').Value = '  +7.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.107'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = $defaultStyleRef.Style
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('D15').Value = '2.712.92'
$ws.Range('E15').Value = '  -0.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '16.58'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = $defaultStyleRef.Style
$ws.Range('E16').Value = '  -2.13%  '
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Value = '2.363.90'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').Value = '43.971.10'
$ws.Range('E19').Value = '  -1.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.96'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = $defaultStyleRef.Style
$ws.Range('E20').Value = '  +7.09%  '
$ws.Range('E21').Value = '  -1.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '77.29'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = $defaultStyleRef.Style
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '259.03'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = $defaultStyleRef.Style
$ws.Range('E23').Value = '  +1.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.99'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = $defaultStyleRef.Style
$ws.Range('E24').Value = '  +24.08%  '
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('E26').Value = '  -3.46%  '
$ws.Range('E27').Value = '  -1.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.81'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = $defaultStyleRef.Style
$ws.Range('E28').Value = '  +3.45%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.24'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = $defaultStyleRef.Style
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.80'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = $defaultStyleRef.Style
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.59'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = $defaultStyleRef.Style
$ws.Range('E31').Value = '  +0.71%  '
$ws.Range('E32').Value = '  -2.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.137'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('D33').Style = $defaultStyleRef.Style
$ws.Range('E33').Value = '  +2.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0762'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = $defaultStyleRef.Style
$ws.Range('E34').Value = '  +1.85%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.58'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = $defaultStyleRef.Style
$ws.Range('E35').Value = '  +6.46%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.24'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = $defaultStyleRef.Style
$ws.Range('E36').Value = '  +0.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.76'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = $defaultStyleRef.Style
$ws.Range('E37').Value = '  -4.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.37'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = $defaultStyleRef.Style
$ws.Range('E38').Value = '  -2.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.35'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = $defaultStyleRef.Style
$ws.Range('E39').Value = '  -3.57%  '
$ws.Range('E40').Value = '  +2.41%  '
$ws.Range('E41').Value = '  +14.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.203'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = $defaultStyleRef.Style
$ws.Range('E42').Value = '  +9.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.02'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = $defaultStyleRef.Style
$ws.Range('E43').Value = '  -5.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.00'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = $defaultStyleRef.Style
$ws.Range('E44').Value = '  +1.57%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.75'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = $defaultStyleRef.Style
$ws.Range('E46').Value = '  +5.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '58.71'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = $defaultStyleRef.Style
$ws.Range('E47').Value = '  +11.44%  '
$ws.Range('E48').Value = '  +5.87%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.25'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = $defaultStyleRef.Style
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '100.89'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = $defaultStyleRef.Style
$ws.Range('E50').Value = '  +2.01%  '
$ws.Range('E51').Value = '  -0.04%  '
